# Applies the "more test cases on activities and items" commit:
#  - updates selections on two existing sheets
#  - appends 7 new worksheets with their data
#  - leaves the final new sheet ("Add_new_activity") as the active/selected tab

$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, [string]$addr, [string]$val) {
    # Force text storage (shared string) even for numeric-looking values like
    # "12" / "333" so they come out as <c t="s"> rather than numeric <v>.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

function Add-SheetAtEnd($workbook, [string]$name) {
    $lastSheet = $workbook.Worksheets.Item($workbook.Worksheets.Count)
    $newSheet = $workbook.Worksheets.Add([Type]::Missing, $lastSheet)
    $newSheet.Name = $name
    return $newSheet
}

# ---------------------------------------------------------------------------
# 1. Existing sheet selection tweaks
# ---------------------------------------------------------------------------

# Add_Partner_Field_Validation: selection was the whole-sheet sqref, now A2
$wsFieldValidation = $wb.Worksheets.Item("Add_Partner_Field_Validation")
$wsFieldValidation.Range("A2").Select()

# update_Partner_Company_Info: was the tab-selected sheet with E2 selected;
# it loses tabSelected (the new Add_new_activity sheet takes it) and its
# selection moves to A2.
$wsCompanyInfo = $wb.Worksheets.Item("update_Partner_Company_Info")
$wsCompanyInfo.Range("A2").Select()

# ---------------------------------------------------------------------------
# 2. New sheet: update_partner_agreeement_info
# ---------------------------------------------------------------------------
$ws18 = Add-SheetAtEnd $wb "update_partner_agreeement_info"
$ws18.Range("A1").Value = "bookingtype"
$ws18.Range("B1").Value = "website"
$ws18.Range("C1").Value = "username"
$ws18.Range("D1").Value = "password"
$ws18.Range("E1").Value = "generalnotes"
$ws18.Range("F1").Value = "runmode"
$ws18.Range("A2").Value = "DIRECT"
$ws18.Range("B2").Value = "google.com"
$ws18.Range("C2").Value = "larry"
$ws18.Range("D2").Value = "rajapaksa"
$ws18.Range("E2").Value = "very good partner to deal with"
$ws18.Range("F2").Value = "Y"
$ws18.Range("B5").Select()

# ---------------------------------------------------------------------------
# 3. New sheet: upload_partner_agreement
# ---------------------------------------------------------------------------
$ws19 = Add-SheetAtEnd $wb "upload_partner_agreement"
$ws19.Range("A1").Value = "file"
$ws19.Range("B1").Value = "documentname"
$ws19.Range("C1").Value = "toastmessage"
$ws19.Range("D1").Value = "runmode"
$ws19.Range("A2").Value = "src/test/resources/images/test.pdf"
$ws19.Range("B2").Value = "Ag13"
$ws19.Range("C2").Value = "The uploaded file will appear in agreements soon"
$ws19.Range("D2").Value = "Y"
$ws19.Range("E9").Select()

# ---------------------------------------------------------------------------
# 4. New sheet: Add_new_activity (ends up the active / tab-selected sheet)
# ---------------------------------------------------------------------------
$ws20 = Add-SheetAtEnd $wb "Add_new_activity"
$ws20.Range("A1").Value = "activityType"
$ws20.Range("B1").Value = "title"
$ws20.Range("C1").Value = "activityimage"
$ws20.Range("D1").Value = "desc"
$ws20.Range("E1").Value = "toastmessage"
$ws20.Range("F1").Value = "runmode"
$ws20.Range("A2").Value = "Flight"
$ws20.Range("B2").Value = "Colombo to Norway"
$ws20.Range("C2").Value = "src/test/resources/images/logo-color.png"
$ws20.Range("D2").Value = "nothing much to say2"
$ws20.Range("E2").Value = "Activity added successfully!"
$ws20.Range("F2").Value = "Y"
$ws20.Range("A3").Value = "Flight"
$ws20.Range("B3").Value = "India to US"
$ws20.Range("C3").Value = "src/test/resources/images/logo-color.png"
$ws20.Range("D3").Value = "Second iteration"
$ws20.Range("E3").Value = "Activity added successfully!"
$ws20.Range("F3").Value = "Y"
$ws20.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. New sheet: update_activity_template
# ---------------------------------------------------------------------------
$ws21 = Add-SheetAtEnd $wb "update_activity_template"
$ws21.Range("A1").Value = "desc"
$ws21.Range("B1").Value = "toastmessage"
$ws21.Range("C1").Value = "runmode"
$ws21.Range("A2").Value = "Updated the description 1"
$ws21.Range("B2").Value = "Activity updated successfully!"
$ws21.Range("C2").Value = "Y"
$ws21.Range("A3").Select()

# ---------------------------------------------------------------------------
# 6. New sheet: add_new_item (values that look numeric are forced to text)
# ---------------------------------------------------------------------------
$ws22 = Add-SheetAtEnd $wb "add_new_item"
Set-TextCell $ws22 "A1" "runmode"
Set-TextCell $ws22 "B1" "itemname"
Set-TextCell $ws22 "C1" "priceoption"
Set-TextCell $ws22 "D1" "op1_upn_price"
Set-TextCell $ws22 "E1" "op1_vat_rate"
Set-TextCell $ws22 "F1" "op1_marketprice"
Set-TextCell $ws22 "G1" "op1_commission"
Set-TextCell $ws22 "H1" "op2_marketpriceadult"
Set-TextCell $ws22 "I1" "op2_marketpricechild"
Set-TextCell $ws22 "J1" "op3_priceroom"
Set-TextCell $ws22 "K1" "desc"
Set-TextCell $ws22 "L1" "toastmessage"

Set-TextCell $ws22 "A2" "Y"
Set-TextCell $ws22 "B2" "arshad_iitem_1_fixed_price"
Set-TextCell $ws22 "C2" "FIXED"
Set-TextCell $ws22 "D2" "567"
Set-TextCell $ws22 "E2" "12"
Set-TextCell $ws22 "F2" "333"
Set-TextCell $ws22 "G2" "3"
Set-TextCell $ws22 "K2" "long desc long desc long desc long desc long desc long desc long desc long desc long desc long desc"
Set-TextCell $ws22 "L2" "Item added successfully!"

Set-TextCell $ws22 "A3" "N"
Set-TextCell $ws22 "B3" "arshad_iitem_2_Person"
Set-TextCell $ws22 "C3" "PERSON"
Set-TextCell $ws22 "E3" "15"
Set-TextCell $ws22 "G3" "7"
Set-TextCell $ws22 "H3" "333"
Set-TextCell $ws22 "I3" "444"
Set-TextCell $ws22 "K3" "Recently UpdatedRecently UpdatedRecently UpdatedRecently Updated Recently Updated Recently UpdatedRecently UpdatedRecently UpdatedRecently UpdatedRecently Updated"
Set-TextCell $ws22 "L3" "Item added successfully!"

Set-TextCell $ws22 "A4" "N"
Set-TextCell $ws22 "B4" "arshad_iitem_3_Room"
Set-TextCell $ws22 "C4" "ROOM"
Set-TextCell $ws22 "E4" "25"
Set-TextCell $ws22 "G4" "9"
Set-TextCell $ws22 "J4" "777"
Set-TextCell $ws22 "K4" "Only My Issues Only My Issues Only My Issues Only My Issues Only My Issues Only My Issues Only My Issues Only My Issues Only My Issues"
Set-TextCell $ws22 "L4" "Item added successfully!"

$ws22.Range("A2").Select()
$ws22.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 7. New sheet: edit_item_template
# ---------------------------------------------------------------------------
$ws23 = Add-SheetAtEnd $wb "edit_item_template"
$ws23.Range("A1").Value = "runmode"
$ws23.Range("A2").Value = "Y"
$ws23.Range("B2").Select()

# ---------------------------------------------------------------------------
# 8. New sheet: add_another_item
# ---------------------------------------------------------------------------
$ws24 = Add-SheetAtEnd $wb "add_another_item"
$ws24.Range("A1").Value = "runmode"
$ws24.Range("B1").Value = "toastmessage"
$ws24.Range("A2").Value = "Y"
$ws24.Range("B2").Value = "Item added successfully!"
$ws24.Range("B2").Select()

# ---------------------------------------------------------------------------
# 9. Final active sheet / selection: Add_new_activity, cell F4
# ---------------------------------------------------------------------------
$ws20.Range("F4").Select()
